$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# extend mustersan : 38.2-38 => 38.2-35
$ws.Range("D8").Value = 35
$ws.Range("C9").Value = 35

# swap C10/D10 values: C10 38.2 (was 39), D10 39 (was 38.2)
$ws.Range("C10").Value = 38.2
$ws.Range("D10").Value = 39

# update selection to D20
$ws.Range("D20").Select()
